# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-12, 14-15 (row 13 unchanged)
$kValues = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 0
    14 = 2
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
